$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look like pure numbers,
# so Excel stores them as literal text (preserving trailing zeros / exact digits)
# instead of silently coercing to a numeric type.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "52.187.96"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "2.833.44"
$ws.Range("E3").Value = "  +3.43%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "355.23"
$ws.Range("E5").Value = "  +6.57%  "
$ws.Range("D6").Value = "113.88"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("E7").Value = "  +2.67%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.611"
$ws.Range("E9").Value = "  +7.52%  "
$ws.Range("D10").Value = "42.17"
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("D11").Value = "0.0851"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "20.29"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "7.82"
$ws.Range("E14").Value = "  +3.74%  "
$ws.Range("D15").Value = "3.260.87"
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("D16").Value = "2.832.49"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("D17").Value = "0.896"
$ws.Range("E17").Value = "  +2.55%  "
$ws.Range("D18").Value = "52.200.57"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").Value = "13.85"
$ws.Range("E19").Value = "  +3.26%  "
$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").Value = "3.18"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "7.31"
$ws.Range("E21").Value = "  +7.31%  "
$ws.Range("D22").Value = "0.0₃0999"
$ws.Range("E22").Value = "  +2.91%  "
$ws.Range("D23").Value = "270.74"
$ws.Range("E23").Value = "  -2.89%  "
$ws.Range("D24").Value = "69.75"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +5.99%  "
$ws.Range("D26").Value = "26.72"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "10.29"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("D32").Value = "33.96"
$ws.Range("E32").Value = "  -2.47%  "
$ws.Range("D33").Value = "5.92"
$ws.Range("E33").Value = "  +6.93%  "
$ws.Range("D34").Value = "0.0445"
$ws.Range("E34").Value = "  +29.32%  "
$ws.Range("D35").Value = "0.0834"
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "2.10"
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "4.90"
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "18.50"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").Value = "3.22"
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("E41").Value = "  +9.44%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "23.60"
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "128.15"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").Value = "3.36"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").Value = "2.046.72"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("E48").Value = "  +3.57%  "
$ws.Range("D49").Value = "0.966"
$ws.Range("E49").Value = "  +11.82%  "
$ws.Range("E50").Value = "  +3.34%  "
$ws.Range("D51").Value = "60.56"
$ws.Range("E51").Value = "  +1.66%  "
